$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "A"
$ws.Range("B4").Value = "PRM3"

$ws.Range("A5").Value = "Loris Lefebvre"
$ws.Range("B5").Value = "PRM4"
$ws.Range("C5").Value = "04h - 12h "
$ws.Range("D5").Value = "A"
$ws.Range("E5").Value = "OUI"
$ws.Range("F5").Value = "OUI"

$ws.Range("F5").Select()
